$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.292.18"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.663.84"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "218.58"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "0.5326"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "20.55"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "4.554"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "1.665.93"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "0.5528"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "65.70"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").Value = "192.39"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "10.16"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "6.058"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "145.07"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").Value = "0.1225"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "7.243"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "16.13"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "1.475"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").Value = "0.05850"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("D32").Value = "3.310"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("D34").Value = "2.823"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").Value = "0.9593"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("D36").Value = "2.432"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "0.5815"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").Value = "0.01612"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "5.881"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "0.8541"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.010"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.047.93"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "104.79"
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("D44").Value = "1.805.84"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "7.933"
$ws.Range("D51").Value = "1.446"
$ws.Range("E51").Value = "  -1.65%  "
